$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task04")

$ws.Range("D9").Value = "https://kimcoder.tistory.com/244?category=911141"
$ws.Range("E9").Value = "security 정리 잘 됨"

$ws.Range("D10").Value = "https://12716.tistory.com/entry/%EC%95%8C%EA%B3%A0%EB%A6%AC%EC%A6%98%EB%B0%B1%EC%A4%80Baekjoon%EB%B0%B1%EC%A4%80-NO10430-Java%EC%9E%90%EB%B0%94%EB%A1%9C-%ED%92%80%EC%96%B4%EB%B3%B4%EA%B8%B0"
$ws.Range("E10").Value = "405 에러"

$ws.Range("D11").Value = "https://codevang.tistory.com/268"
$ws.Range("E11").Value = "security 정리 잘 됨 222"

$ws.Range("D12").Value = "https://taesan94.tistory.com/124"

$ws.Range("D10").Select()
